$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

$ws.Range("C8").Value = "El equipo ha completado los productos especificados. El equipo ha acumulado toda la información y ha completado todas las formas requeridas."
$ws.Range("C6").Value = "Se ha implementado la funcionalidad #3 del producto."
$ws.Range("C2").Value = "El equipo ha completado un ciclo preeviamente."
$ws.Range("C3").Value = "El equipo ha completado un ciclo preeviamente. Cada miembro del equipo ha leído el capítulo correspondiente a su rol."
$ws.Range("C4").Value = "Se ha completado la reunión de equipo para analizar la versión final del documento de requerimientos. Se ha creado el esquema del documento de arquitectura."
$ws.Range("C5").Value = ""
$ws.Range("C7").Value = ""

$ws.Range("D4").Select()
